# Add 2022-Q4 data:
#  - insert a new worksheet "2022-Q4" right before "2022-Q3"
#  - populate it with the Q4 fund-holdings table
#  - add a new row to the "总计" (summary) sheet for 2022-Q4 and
#    shift the existing rows down

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q4" worksheet before "2022-Q3"
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($q3)
$q4.Name = "2022-Q4"

# Header row
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $cell = $q4.Cells.Item(1, 2 + $c)
    $cell.Value = $headers[$c]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# Data rows: index, code, name, scale, position, positionPct, marketValue, rank
$rows = @(
    @(0,  "000601", "华宝创新优选混合",               "10.34", "87.48", "1.99", "0.2058", 10),
    @(1,  "006081", "海富通电子信息传媒产业股票A",     "5.43",  "92.10", "3.47", "0.1884", 6),
    @(2,  "006265", "红土创新新科技股票",               "2.90",  "93.37", "5.01", "0.1453", 7),
    @(3,  "006080", "海富通电子信息传媒产业股票C",     "3.10",  "92.10", "3.47", "0.1076", 6),
    @(4,  "004818", "国寿安保目标策略灵活配置混合A",   "2.76",  "59.92", "3.36", "0.0927", 3),
    @(5,  "006449", "浙商汇金量化精选灵活配置混合",     "2.47",  "59.45", "2.77", "0.0684", 9),
    @(6,  "005903", "泰达宏利绩优增长灵活配置混合A",   "1.85",  "87.04", "2.46", "0.0455", 10),
    @(7,  "004819", "国寿安保目标策略灵活配置混合C",   "1.30",  "59.92", "3.36", "0.0437", 3),
    @(8,  "015576", "泰达宏利绩优增长灵活配置混合C",   "1.70",  "87.04", "2.46", "0.0418", 10),
    @(9,  "016013", "南方碳中和股票A",                   "0.42",  "84.75", "4.06", "0.0171", 4),
    @(10, "001659", "富安达新动力灵活配置混合",         "0.11",  "86.09", "2.98", "0.0033", 9),
    @(11, "016014", "南方碳中和股票C",                   "0.01",  "84.75", "4.06", "0.0004", 4)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]

    $idxCell = $q4.Cells.Item($r, 1)
    $idxCell.Value = $row[0]
    $idxCell.Font.Bold = $true
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160
    $idxCell.Borders.LineStyle = 1

    # text-valued columns B..G (fund code, name, scale, position, pct, market value)
    for ($col = 2; $col -le 7; $col++) {
        $cell = $q4.Cells.Item($r, $col)
        $cell.NumberFormat = "@"
        $cell.Value = [string]$row[$col - 1]
    }

    # rank column H is numeric
    $q4.Cells.Item($r, 8).Value = $row[7]
}

# ---------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: insert a row for 2022-Q4
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

$a2 = $summary.Cells.Item(2, 1)
$a2.Value = 0
$a2.Font.Bold = $true
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160
$a2.Borders.LineStyle = 1

$summary.Cells.Item(2, 2).Value = "2022-Q4"
$summary.Cells.Item(2, 3).Value = 12
$summary.Cells.Item(2, 4).Value = 0.96

# renumber the index column (A) for the rows that shifted down
$summary.Cells.Item(3, 1).Value = 1
$summary.Cells.Item(4, 1).Value = 2
$summary.Cells.Item(5, 1).Value = 3
